$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

# New asset rows: str_batchName and bol_containError
$ws.Range("A4").Value = "str_batchName"
$ws.Range("B4").Value = "str_batchName"
$ws.Range("C4").Value = "logement"
$ws.Range("D4").Value = "nom du batch"

$ws.Range("A5").Value = "bol_containError"
$ws.Range("B5").Value = "bol_containError"
$ws.Range("C5").Value = "logement"
$ws.Range("D5").Value = "presence erreur dans le batch"

# Update the active selection on the Assets sheet to A5 (matches author's editing position)
$ws.Activate()
$ws.Range("A5").Select()
